$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(16, 8).Value = 711.6
$ws.Cells.Item(16, 9).Value = 711.6
$ws.Cells.Item(16, 11).Value = 711.6
$ws.Cells.Item(16, 13).Value = -481.6

$ws.Cells.Item(80, 8).Value = 403.73334
$ws.Cells.Item(80, 9).Value = 446.75
$ws.Cells.Item(80, 10).Value = 231.66667
$ws.Cells.Item(80, 11).Value = 1340.25
$ws.Cells.Item(80, 12).Value = 695.00001
$ws.Cells.Item(80, 13).Value = -342.25
$ws.Cells.Item(80, 14).Value = -2691.00001

$ws.Cells.Item(83, 8).Value = 403.73334
$ws.Cells.Item(83, 9).Value = 446.75
$ws.Cells.Item(83, 10).Value = 231.66667
$ws.Cells.Item(83, 11).Value = 4020.75
$ws.Cells.Item(83, 12).Value = 2085.00003
$ws.Cells.Item(83, 13).Value = 971.25
$ws.Cells.Item(83, 14).Value = -12069.00003

$ws.Cells.Item(92, 8).Value = 5051367.5
$ws.Cells.Item(92, 9).Value = 6173790
$ws.Cells.Item(92, 10).Value = 465
$ws.Cells.Item(92, 11).Value = 6173790
$ws.Cells.Item(92, 12).Value = 465
$ws.Cells.Item(92, 13).Value = -6172542
$ws.Cells.Item(92, 14).Value = -2961

$ws.Cells.Item(98, 8).Value = 65790132
$ws.Cells.Item(98, 9).Value = 69445020
$ws.Cells.Item(98, 10).Value = 2000
$ws.Cells.Item(98, 11).Value = 69445020
$ws.Cells.Item(98, 12).Value = 2000
$ws.Cells.Item(98, 13).Value = -69443522
$ws.Cells.Item(98, 14).Value = -4996

$ws.Cells.Item(122, 8).Value = 65790132
$ws.Cells.Item(122, 9).Value = 69445020
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 208335060
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -208332610
$ws.Cells.Item(122, 14).Value = -10900

$ws.Cells.Item(135, 8).Value = 1117.2354
$ws.Cells.Item(135, 9).Value = 1068.3125
$ws.Cells.Item(135, 10).Value = 1900
$ws.Cells.Item(135, 11).Value = 9614.8125
$ws.Cells.Item(135, 12).Value = 17100
$ws.Cells.Item(135, 13).Value = -7079.8125
$ws.Cells.Item(135, 14).Value = -22170

$ws.Cells.Item(137, 8).Value = 1414.8286
$ws.Cells.Item(137, 9).Value = 1294.4546
$ws.Cells.Item(137, 10).Value = 1618.5385
$ws.Cells.Item(137, 11).Value = 3883.3638
$ws.Cells.Item(137, 12).Value = 4855.6155
$ws.Cells.Item(137, 13).Value = -1333.3638
$ws.Cells.Item(137, 14).Value = -9955.6155

$ws.Cells.Item(138, 8).Value = 1334.3677
$ws.Cells.Item(138, 9).Value = 743.9245
$ws.Cells.Item(138, 10).Value = 3420.6
$ws.Cells.Item(138, 11).Value = 2231.7735
$ws.Cells.Item(138, 12).Value = 10261.8
$ws.Cells.Item(138, 13).Value = 2908.2265
$ws.Cells.Item(138, 14).Value = -20541.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1148.0714
$ws.Cells.Item(61, 9).Value = 1102.4474
$ws.Cells.Item(61, 10).Value = 1244.3889
$ws.Cells.Item(61, 11).Value = 1102.4474
$ws.Cells.Item(61, 12).Value = 1244.3889
$ws.Cells.Item(61, 13).Value = -890.4474
$ws.Cells.Item(61, 14).Value = -1668.3889

$ws.Cells.Item(74, 8).Value = 1375.7894
$ws.Cells.Item(74, 9).Value = 1423.4706
$ws.Cells.Item(74, 10).Value = 970.5
$ws.Cells.Item(74, 11).Value = 1423.4706
$ws.Cells.Item(74, 12).Value = 970.5
$ws.Cells.Item(74, 13).Value = -549.4706000000001
$ws.Cells.Item(74, 14).Value = -2718.5

$ws.Cells.Item(77, 8).Value = 1375.7894
$ws.Cells.Item(77, 9).Value = 1423.4706
$ws.Cells.Item(77, 10).Value = 970.5
$ws.Cells.Item(77, 11).Value = 7117.353000000001
$ws.Cells.Item(77, 12).Value = 4852.5
$ws.Cells.Item(77, 13).Value = -2749.353000000001
$ws.Cells.Item(77, 14).Value = -13588.5

$ws.Cells.Item(113, 8).Value = 49755.145
$ws.Cells.Item(113, 10).Value = 49755.145
$ws.Cells.Item(113, 12).Value = 49755.145
$ws.Cells.Item(113, 14).Value = -58433.145

$ws.Cells.Item(122, 8).Value = 1178.8
$ws.Cells.Item(122, 9).Value = 973.5
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 2920.5
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -470.5
$ws.Cells.Item(122, 14).Value = -10900

$ws.Cells.Item(136, 8).Value = 1148.0714
$ws.Cells.Item(136, 9).Value = 1102.4474
$ws.Cells.Item(136, 10).Value = 1244.3889
$ws.Cells.Item(136, 11).Value = 3307.3422
$ws.Cells.Item(136, 12).Value = 3733.1667
$ws.Cells.Item(136, 13).Value = -757.3422
$ws.Cells.Item(136, 14).Value = -8833.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1663431.5
$ws.Cells.Item(86, 9).Value = 2719.8572
$ws.Cells.Item(86, 10).Value = 3324143.2
$ws.Cells.Item(86, 11).Value = 2719.8572
$ws.Cells.Item(86, 12).Value = 3324143.2
$ws.Cells.Item(86, 13).Value = -1596.8572
$ws.Cells.Item(86, 14).Value = -3326389.2

$ws.Cells.Item(89, 8).Value = 1663431.5
$ws.Cells.Item(89, 9).Value = 2719.8572
$ws.Cells.Item(89, 10).Value = 3324143.2
$ws.Cells.Item(89, 11).Value = 13599.286
$ws.Cells.Item(89, 12).Value = 16620716
$ws.Cells.Item(89, 13).Value = -7983.286
$ws.Cells.Item(89, 14).Value = -16631948

$ws.Cells.Item(94, 8).Value = 668.3158
$ws.Cells.Item(94, 9).Value = 586.5333000000001
$ws.Cells.Item(94, 10).Value = 975
$ws.Cells.Item(94, 11).Value = 586.5333000000001
$ws.Cells.Item(94, 12).Value = 975
$ws.Cells.Item(94, 13).Value = -135.5333000000001
$ws.Cells.Item(94, 14).Value = -1877

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1251.1578
$ws.Cells.Item(31, 9).Value = 994.12726
$ws.Cells.Item(31, 10).Value = 1924.3334
$ws.Cells.Item(31, 11).Value = 994.12726
$ws.Cells.Item(31, 12).Value = 1924.3334
$ws.Cells.Item(31, 13).Value = -699.12726
$ws.Cells.Item(31, 14).Value = -2514.3334

$ws.Cells.Item(34, 8).Value = 1251.1578
$ws.Cells.Item(34, 9).Value = 994.12726
$ws.Cells.Item(34, 10).Value = 1924.3334
$ws.Cells.Item(34, 11).Value = 994.12726
$ws.Cells.Item(34, 12).Value = 1924.3334
$ws.Cells.Item(34, 13).Value = -792.12726
$ws.Cells.Item(34, 14).Value = -2328.3334

$ws.Cells.Item(81, 8).Value = 62333.332
$ws.Cells.Item(81, 10).Value = 62333.332
$ws.Cells.Item(81, 12).Value = 62333.332
$ws.Cells.Item(81, 14).Value = -64329.332

$ws.Cells.Item(84, 8).Value = 62333.332
$ws.Cells.Item(84, 10).Value = 62333.332
$ws.Cells.Item(84, 12).Value = 186999.996
$ws.Cells.Item(84, 14).Value = -196983.996

$ws.Cells.Item(99, 8).Value = 166668850
$ws.Cells.Item(99, 9).Value = 200001730
$ws.Cells.Item(99, 10).Value = 4500
$ws.Cells.Item(99, 11).Value = 200001730
$ws.Cells.Item(99, 12).Value = 4500
$ws.Cells.Item(99, 13).Value = -200000232
$ws.Cells.Item(99, 14).Value = -7496

$ws.Cells.Item(126, 8).Value = 166668850
$ws.Cells.Item(126, 9).Value = 200001730
$ws.Cells.Item(126, 10).Value = 4500
$ws.Cells.Item(126, 11).Value = 600005190
$ws.Cells.Item(126, 12).Value = 13500
$ws.Cells.Item(126, 13).Value = -600002720
$ws.Cells.Item(126, 14).Value = -18440

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 2300.52
$ws.Cells.Item(131, 9).Value = 2218.3333
$ws.Cells.Item(131, 10).Value = 2305.7659
$ws.Cells.Item(131, 11).Value = 6654.999899999999
$ws.Cells.Item(131, 12).Value = 6917.297699999999
$ws.Cells.Item(131, 13).Value = -1614.999899999999
$ws.Cells.Item(131, 14).Value = -16997.2977

$ws.Cells.Item(132, 8).Value = 12198628
$ws.Cells.Item(132, 9).Value = 713.1539
$ws.Cells.Item(132, 10).Value = 17861944
$ws.Cells.Item(132, 11).Value = 6418.3851
$ws.Cells.Item(132, 12).Value = 160757496
$ws.Cells.Item(132, 13).Value = -3888.3851
$ws.Cells.Item(132, 14).Value = -160762556

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 2101382.5
$ws.Cells.Item(12, 9).Value = 3001331.5
$ws.Cells.Item(12, 10).Value = 1501.3334
$ws.Cells.Item(12, 11).Value = 3001331.5
$ws.Cells.Item(12, 12).Value = 1501.3334
$ws.Cells.Item(12, 13).Value = -3001191.5
$ws.Cells.Item(12, 14).Value = -1781.3334

$ws.Cells.Item(80, 8).Value = 5002936
$ws.Cells.Item(80, 9).Value = 3270
$ws.Cells.Item(80, 10).Value = 25001600
$ws.Cells.Item(80, 11).Value = 3270
$ws.Cells.Item(80, 12).Value = 25001600
$ws.Cells.Item(80, 13).Value = -2272
$ws.Cells.Item(80, 14).Value = -25003596

$ws.Cells.Item(83, 8).Value = 5002936
$ws.Cells.Item(83, 9).Value = 3270
$ws.Cells.Item(83, 10).Value = 25001600
$ws.Cells.Item(83, 11).Value = 16350
$ws.Cells.Item(83, 12).Value = 125008000
$ws.Cells.Item(83, 13).Value = -11358
$ws.Cells.Item(83, 14).Value = -125017984

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1574.6428
$ws.Cells.Item(7, 9).Value = 1460.5454
$ws.Cells.Item(7, 10).Value = 1993
$ws.Cells.Item(7, 11).Value = 1460.5454
$ws.Cells.Item(7, 12).Value = 1993
$ws.Cells.Item(7, 13).Value = -1348.5454
$ws.Cells.Item(7, 14).Value = -2217

$ws.Cells.Item(16, 8).Value = 448.5
$ws.Cells.Item(16, 9).Value = 407.27274
$ws.Cells.Item(16, 10).Value = 599.6667
$ws.Cells.Item(16, 11).Value = 407.27274
$ws.Cells.Item(16, 12).Value = 599.6667
$ws.Cells.Item(16, 13).Value = -237.27274
$ws.Cells.Item(16, 14).Value = -939.6667

$ws.Cells.Item(40, 8).Value = 50002700
$ws.Cells.Item(40, 9).Value = 3666.6667
$ws.Cells.Item(40, 10).Value = 125001250
$ws.Cells.Item(40, 11).Value = 3666.6667
$ws.Cells.Item(40, 12).Value = 125001250
$ws.Cells.Item(40, 13).Value = -3530.6667
$ws.Cells.Item(40, 14).Value = -125001522

$ws.Cells.Item(55, 8).Value = 5482.2104
$ws.Cells.Item(55, 9).Value = 14551.571
$ws.Cells.Item(55, 11).Value = 14551.571
$ws.Cells.Item(55, 13).Value = -14378.571

$ws.Cells.Item(126, 8).Value = 1574.6428
$ws.Cells.Item(126, 9).Value = 1460.5454
$ws.Cells.Item(126, 10).Value = 1993
$ws.Cells.Item(126, 11).Value = 4381.6362
$ws.Cells.Item(126, 12).Value = 5979
$ws.Cells.Item(126, 13).Value = -1911.6362
$ws.Cells.Item(126, 14).Value = -10919

$ws.Cells.Item(136, 8).Value = 31748288
$ws.Cells.Item(136, 9).Value = 5293675.5
$ws.Cells.Item(136, 11).Value = 15881026.5
$ws.Cells.Item(136, 13).Value = -15878476.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 10000
$ws.Cells.Item(20, 9).Value = 10000
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 10000
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -9760
$ws.Range("N20").ClearContents()

$ws.Cells.Item(126, 8).Value = 1311.3572
$ws.Cells.Item(126, 9).Value = 993.1667
$ws.Cells.Item(126, 10).Value = 1550
$ws.Cells.Item(126, 11).Value = 2979.5001
$ws.Cells.Item(126, 12).Value = 4650
$ws.Cells.Item(126, 13).Value = -509.5001000000002
$ws.Cells.Item(126, 14).Value = -9590

$ws.Cells.Item(132, 8).Value = 5695958.5
$ws.Cells.Item(132, 9).Value = 16550.637
$ws.Cells.Item(132, 10).Value = 22734182
$ws.Cells.Item(132, 11).Value = 49651.91099999999
$ws.Cells.Item(132, 12).Value = 68202546
$ws.Cells.Item(132, 13).Value = -47121.91099999999
$ws.Cells.Item(132, 14).Value = -68207606

$ws.Cells.Item(136, 8).Value = 9095313
$ws.Cells.Item(136, 9).Value = 13894569
$ws.Cells.Item(136, 10).Value = 1987.0526
$ws.Cells.Item(136, 11).Value = 41683707
$ws.Cells.Item(136, 12).Value = 5961.1578
$ws.Cells.Item(136, 13).Value = -41681157
$ws.Cells.Item(136, 14).Value = -11061.1578
